$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell
$ws.Range("F1").Value = "DateNaissance"

# F2: first real date value, then apply the built-in short-date format (numFmtId 14)
$ws.Range("F2").Value = (Get-Date -Year 2002 -Month 11 -Day 28 -Hour 0 -Minute 0 -Second 0)
$ws.Range("F2").NumberFormat = "mm-dd-yy"
$ws.Range("F2").Copy()

# Reuse the exact same style for the other valid-date cells by pasting formats only,
# then assigning their values (this avoids Excel re-creating a duplicate autoformat style).
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").Value = (Get-Date -Year 2002 -Month 2 -Day 16 -Hour 0 -Minute 0 -Second 0)

$ws.Range("F7").PasteSpecial(-4122)
$ws.Range("F7").Value = (Get-Date -Year 2009 -Month 11 -Day 12 -Hour 0 -Minute 0 -Second 0)

$ws.Range("F8").PasteSpecial(-4122)
$ws.Range("F8").Value = (Get-Date -Year 2002 -Month 12 -Day 13 -Hour 0 -Minute 0 -Second 0)

$ws.Range("F9").PasteSpecial(-4122)
$ws.Range("F9").Value = (Get-Date -Year 2000 -Month 4 -Day 23 -Hour 0 -Minute 0 -Second 0)

# Invalid "dates" stored as plain text
$ws.Range("F3").Value = "244/12/2003"
$ws.Range("F4").Value = "12/38-2004"
$ws.Range("F6").Value = "51/13/2001"

# New data rows 8 and 9
$ws.Range("A8").Value = 456
$ws.Range("B8").Value = "hfg"
$ws.Range("C8").Value = "dfdd"
$ws.Range("D8").Value = 12333222
$ws.Range("E8").Value = "dfdd@gmail.com"

$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Fahim"
$ws.Range("C9").Value = "Malika"
$ws.Range("D9").Value = 1234567
$ws.Range("E9").Value = "Malika@gmail.com"

$ws.Hyperlinks.Add($ws.Range("E8"), "mailto:dfdd@gmail.com")
$ws.Hyperlinks.Add($ws.Range("E9"), "mailto:Malika@gmail.com")

$ws.Columns.Item(6).ColumnWidth = 13.42578125

$ws.Range("F10").Select()
